$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 433.77777
$ws.Range("I53").Value = 160.4
$ws.Range("K53").Value = 160.4
$ws.Range("M53").Value = 476.6
$ws.Range("H112").Value = 1432.8572
$ws.Range("I112").Value = 825
$ws.Range("J112").Value = 1534.1666
$ws.Range("K112").Value = 2475
$ws.Range("L112").Value = 4602.4998
$ws.Range("M112").Value = -1367
$ws.Range("N112").Value = -6818.4998
$ws.Range("H129").Value = 868.75
$ws.Range("J129").Value = 1816.6666
$ws.Range("L129").Value = 5449.9998
$ws.Range("N129").Value = -15449.9998
$ws.Range("H132").Value = 1003170.6
$ws.Range("I132").Value = 4100.7144
$ws.Range("K132").Value = 12302.1432
$ws.Range("M132").Value = -9772.143199999999
$ws.Range("H137").Value = 2798.111
$ws.Range("I137").Value = 3264
$ws.Range("J137").Value = 1866.3334
$ws.Range("K137").Value = 9792
$ws.Range("L137").Value = 5599.0002
$ws.Range("M137").Value = -7242
$ws.Range("N137").Value = -10699.0002
$ws.Range("H138").Value = 2878.9
$ws.Range("I138").Value = 1049
$ws.Range("J138").Value = 3201.8235
$ws.Range("K138").Value = 3147
$ws.Range("L138").Value = 9605.470499999999
$ws.Range("M138").Value = 1993
$ws.Range("N138").Value = -19885.4705
$ws.Range("H139").Value = 32468.9
$ws.Range("J139").Value = 33775.555
$ws.Range("L139").Value = 33775.555
$ws.Range("N139").Value = -44055.555

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1939.8125
$ws.Range("I2").Value = 1977.909
$ws.Range("J2").Value = 1856
$ws.Range("K2").Value = 1977.909
$ws.Range("L2").Value = 1856
$ws.Range("M2").Value = -1864.909
$ws.Range("N2").Value = -2082
$ws.Range("H32").Value = 6444.4463
$ws.Range("I32").Value = 5943.5
$ws.Range("K32").Value = 5943.5
$ws.Range("M32").Value = -5656.5
$ws.Range("H45").Value = 3176.5
$ws.Range("I45").Value = 2585.3333
$ws.Range("K45").Value = 2585.3333
$ws.Range("M45").Value = -2208.3333
$ws.Range("H61").Value = 4802.8
$ws.Range("I61").Value = 4250
$ws.Range("K61").Value = 4250
$ws.Range("M61").Value = -4038
$ws.Range("H102").Value = 4215.7144
$ws.Range("I102").Value = 4215.7144
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 4215.7144
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -2593.7144
$ws.Range("N102").ClearContents()
$ws.Range("H116").Value = 1939.8125
$ws.Range("I116").Value = 1977.909
$ws.Range("J116").Value = 1856
$ws.Range("K116").Value = 1977.909
$ws.Range("L116").Value = 1856
$ws.Range("M116").Value = 316.0909999999999
$ws.Range("N116").Value = -6444
$ws.Range("H136").Value = 4802.8
$ws.Range("I136").Value = 4250
$ws.Range("K136").Value = 12750
$ws.Range("M136").Value = -10200

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1939.8125
$ws.Range("I3").Value = 1977.909
$ws.Range("J3").Value = 1856
$ws.Range("K3").Value = 1977.909
$ws.Range("L3").Value = 1856
$ws.Range("M3").Value = -1863.909
$ws.Range("N3").Value = -2084
$ws.Range("H134").Value = 85788.75
$ws.Range("I134").Value = 2678.6365
$ws.Range("K134").Value = 8035.9095
$ws.Range("M134").Value = -5500.9095

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2470.3333
$ws.Range("I16").Value = 2011
$ws.Range("J16").Value = 2700
$ws.Range("K16").Value = 2011
$ws.Range("L16").Value = 2700
$ws.Range("M16").Value = -1724
$ws.Range("N16").Value = -3274
$ws.Range("H31").Value = 3826.2068
$ws.Range("I31").Value = 1327.0588
$ws.Range("J31").Value = 7366.6665
$ws.Range("K31").Value = 1327.0588
$ws.Range("L31").Value = 7366.6665
$ws.Range("M31").Value = -1032.0588
$ws.Range("N31").Value = -7956.6665
$ws.Range("H34").Value = 3826.2068
$ws.Range("I34").Value = 1327.0588
$ws.Range("J34").Value = 7366.6665
$ws.Range("K34").Value = 1327.0588
$ws.Range("L34").Value = 7366.6665
$ws.Range("M34").Value = -1125.0588
$ws.Range("N34").Value = -7770.6665
$ws.Range("H62").Value = 9464.706
$ws.Range("I62").Value = 10300
$ws.Range("J62").Value = 8722.223
$ws.Range("K62").Value = 10300
$ws.Range("L62").Value = 8722.223
$ws.Range("M62").Value = -9676
$ws.Range("N62").Value = -9970.223
$ws.Range("H65").Value = 9464.706
$ws.Range("I65").Value = 10300
$ws.Range("J65").Value = 8722.223
$ws.Range("K65").Value = 51500
$ws.Range("L65").Value = 43611.115
$ws.Range("M65").Value = -48380
$ws.Range("N65").Value = -49851.115
$ws.Range("H105").Value = 2964.5
$ws.Range("I105").Value = 2289.3333
$ws.Range("J105").Value = 4990
$ws.Range("K105").Value = 2289.3333
$ws.Range("L105").Value = 4990
$ws.Range("M105").Value = -542.3332999999998
$ws.Range("N105").Value = -8484
$ws.Range("H113").Value = 2470.3333
$ws.Range("I113").Value = 2011
$ws.Range("J113").Value = 2700
$ws.Range("K113").Value = 2011
$ws.Range("L113").Value = 2700
$ws.Range("M113").Value = 159
$ws.Range("N113").Value = -7040

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1111845.5
$ws.Range("I122").Value = 650
$ws.Range("K122").Value = 5850
$ws.Range("M122").Value = -3400
$ws.Range("H131").Value = 87650.836
$ws.Range("I131").Value = 250357.5
$ws.Range("J131").Value = 6297.5
$ws.Range("K131").Value = 751072.5
$ws.Range("L131").Value = 18892.5
$ws.Range("M131").Value = -746032.5
$ws.Range("N131").Value = -28972.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 24835000
$ws.Range("I14").Value = 24835000
$ws.Range("K14").Value = 24835000
$ws.Range("M14").Value = -24834832
$ws.Range("H80").Value = 3430.7896
$ws.Range("I80").Value = 3480.3125
$ws.Range("J80").Value = 3166.6667
$ws.Range("K80").Value = 3480.3125
$ws.Range("L80").Value = 3166.6667
$ws.Range("M80").Value = -2482.3125
$ws.Range("N80").Value = -5162.6667
$ws.Range("H83").Value = 3430.7896
$ws.Range("I83").Value = 3480.3125
$ws.Range("J83").Value = 3166.6667
$ws.Range("K83").Value = 17401.5625
$ws.Range("L83").Value = 15833.3335
$ws.Range("M83").Value = -12409.5625
$ws.Range("N83").Value = -25817.3335
$ws.Range("H139").Value = 27145
$ws.Range("J139").Value = 27145
$ws.Range("L139").Value = 27145
$ws.Range("N139").Value = -37425

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 2000
$ws.Range("J21").Value = 2000
$ws.Range("L21").Value = 2000
$ws.Range("M21").Value = -2348
$ws.Range("H100").Value = 37038372
$ws.Range("I100").Value = 37038372
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 37038372
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -37037831
$ws.Range("N100").ClearContents()
$ws.Range("H138").Value = 25000
$ws.Range("J138").Value = 25000
$ws.Range("L138").Value = 25000
$ws.Range("N138").Value = -35280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4969.154
$ws.Range("I62").Value = 4787.375
$ws.Range("J62").Value = 5260
$ws.Range("K62").Value = 4787.375
$ws.Range("L62").Value = 5260
$ws.Range("M62").Value = -4163.375
$ws.Range("N62").Value = -6508
$ws.Range("H65").Value = 4969.154
$ws.Range("I65").Value = 4787.375
$ws.Range("J65").Value = 5260
$ws.Range("K65").Value = 23936.875
$ws.Range("L65").Value = 26300
$ws.Range("M65").Value = -20816.875
$ws.Range("N65").Value = -32540
$ws.Range("H122").Value = 1447.5454
$ws.Range("I122").Value = 1355.091
$ws.Range("J122").Value = 1540
$ws.Range("K122").Value = 4065.273
$ws.Range("L122").Value = 4620
$ws.Range("M122").Value = -1615.273
$ws.Range("N122").Value = -9520
$ws.Range("H136").Value = 2503.0605
$ws.Range("I136").Value = 2699.8096
$ws.Range("J136").Value = 2158.75
$ws.Range("K136").Value = 8099.4288
$ws.Range("L136").Value = 6476.25
$ws.Range("M136").Value = -5549.4288
$ws.Range("N136").Value = -11576.25
$ws.Range("H138").Value = 39550
$ws.Range("J138").Value = 39550
$ws.Range("L138").Value = 39550
$ws.Range("N138").Value = -49830
